$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 56.9
$ws.Range("N2").Value = 51.53902399942638
$ws.Range("K3").Value = 53.5
$ws.Range("N3").Value = 51.53902399942638
